$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Walk the data rows top-to-bottom, left-to-right (matching column order
# B, D, F) so newly introduced shared strings are registered in the same
# order Excel would encounter them when rebuilding the table on save.
for ($r = 2; $r -le 43; $r++) {
    # harvester: H.BROWN -> S.GISH
    $ws.Cells.Item($r, 2).Value = "S.GISH"

    # experimentDesign stays "90minuteInduction" for every row
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"

    # strain: only rows 2 and 5 keep a value now (typo fixed on row 2,
    # TDY1974 retained on row 5); every other row's strain is cleared.
    if ($r -eq 2) {
        $ws.Cells.Item($r, 6).Value = "KN99alpha"
    } elseif ($r -eq 5) {
        $ws.Cells.Item($r, 6).Value = "TDY1974"
    } else {
        $ws.Cells.Item($r, 6).ClearContents()
    }
}

# Match the saved selection state (single active cell F8).
$ws.Range("F8").Select()
